$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) store numeric-looking values as text
# (matching the source data feed formatting, e.g. "568.39", "19.50", "0.0240",
# thousand-dot-separated prices like "60.703.60", and padded percents like
# "  +2.82%  "). Mark the range as Text before writing so Excel does not
# auto-convert these strings into numeric values and lose formatting
# (trailing zeros, multiple dot separators, subscript digits, etc.).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "60.703.60"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "2.602.15"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "568.39"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "142.27"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("D9").Value = "2.625.34"
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("D10").Value = "6.49"
$ws.Range("E10").Value = "  -2.89%  "
$ws.Range("E11").Value = "  +2.67%  "
$ws.Range("E12").Value = "  -3.91%  "
$ws.Range("E13").Value = "  +7.27%  "
$ws.Range("D14").Value = "3.064.10"
$ws.Range("E14").Value = "  +1.25%  "
$ws.Range("D15").Value = "60.766.24"
$ws.Range("E15").Value = "  +2.80%  "
$ws.Range("E16").Value = "  +5.19%  "
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("D18").Value = "2.613.50"
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("D20").Value = "11.18"
$ws.Range("E20").Value = "  +8.97%  "
$ws.Range("D21").Value = "348.61"
$ws.Range("E21").Value = "  +3.80%  "
$ws.Range("D22").Value = "7.07"
$ws.Range("E22").Value = "  +12.50%  "
$ws.Range("D24").Value = "0.525"
$ws.Range("E24").Value = "  +13.59%  "
$ws.Range("D25").Value = "64.24"
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("D26").Value = "0.996"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("D27").Value = "0.159"
$ws.Range("E27").Value = "  -1.31%  "
$ws.Range("E28").Value = "  +5.64%  "
$ws.Range("D29").Value = "0.0₃0790"
$ws.Range("E29").Value = "  +1.34%  "
$ws.Range("E30").Value = "  +8.91%  "
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").Value = "6.36"
$ws.Range("E32").Value = "  +4.93%  "
$ws.Range("D33").Value = "162.44"
$ws.Range("E33").Value = "  +1.45%  "
$ws.Range("D34").Value = "19.50"
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("D35").Value = "4.23"
$ws.Range("E35").Value = "  +5.54%  "
$ws.Range("D36").Value = "0.955"
$ws.Range("E36").Value = "  +9.46%  "
$ws.Range("D37").Value = "1.21"
$ws.Range("E37").Value = "  +3.43%  "
$ws.Range("E38").Value = "  +7.56%  "
$ws.Range("D39").Value = "37.73"
$ws.Range("E39").Value = "  +0.83%  "
$ws.Range("D40").Value = "0.850"
$ws.Range("E40").Value = "  -2.80%  "
$ws.Range("E41").Value = "  +3.55%  "
$ws.Range("D42").Value = "295.97"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").Value = "140.10"
$ws.Range("E43").Value = "  +7.07%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("E45").Value = "  +0.87%  "
$ws.Range("E46").Value = "  +1.90%  "
$ws.Range("E47").Value = "  +1.77%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "19.50"
$ws.Range("E48").Value = "  +1.98%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "0.0240"
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").Value = "4.87"
$ws.Range("E51").Value = "  +8.29%  "
